$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force literal text storage so numeric-looking strings
    # (e.g. "59.318.06", "0.0995") are not coerced into numbers,
    # then restore the default style so no stray formatting sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '59.318.06'
Set-TextValue $ws.Range("E2") '  +4.95%  '

Set-TextValue $ws.Range("D3") '3.343.54'
Set-TextValue $ws.Range("E3") '  +2.54%  '

Set-TextValue $ws.Range("E4") '  +0.06%  '

Set-TextValue $ws.Range("D5") '410.50'
Set-TextValue $ws.Range("E5") '  +2.72%  '

Set-TextValue $ws.Range("D6") '114.40'
Set-TextValue $ws.Range("E6") '  +2.65%  '

Set-TextValue $ws.Range("D7") '0.588'
Set-TextValue $ws.Range("E7") '  +5.28%  '

Set-TextValue $ws.Range("E8") '  +0.03%  '

Set-TextValue $ws.Range("D9") '0.643'
Set-TextValue $ws.Range("E9") '  +3.39%  '

Set-TextValue $ws.Range("D10") '40.46'
Set-TextValue $ws.Range("E10") '  +2.27%  '

Set-TextValue $ws.Range("D11") '0.0995'
Set-TextValue $ws.Range("E11") '  +4.28%  '

Set-TextValue $ws.Range("D12") '0.143'
Set-TextValue $ws.Range("E12") '  +1.48%  '

Set-TextValue $ws.Range("D13") '3.879.73'
Set-TextValue $ws.Range("E13") '  +2.87%  '

Set-TextValue $ws.Range("D14") '8.55'
Set-TextValue $ws.Range("E14") '  +5.19%  '

Set-TextValue $ws.Range("D15") '19.56'
Set-TextValue $ws.Range("E15") '  +1.37%  '

Set-TextValue $ws.Range("D16") '3.351.52'
Set-TextValue $ws.Range("E16") '  +2.79%  '

Set-TextValue $ws.Range("D17") '1.05'
Set-TextValue $ws.Range("E17") '  +0.15%  '

Set-TextValue $ws.Range("D18") '59.135.20'
Set-TextValue $ws.Range("E18") '  +4.78%  '

Set-TextValue $ws.Range("D19") '10.86'
Set-TextValue $ws.Range("E19") '  -2.01%  '

Set-TextValue $ws.Range("D20") '3.37'
Set-TextValue $ws.Range("E20") '  +1.21%  '

Set-TextValue $ws.Range("D21") '0.0000112'
Set-TextValue $ws.Range("E21") '  +7.25%  '

Set-TextValue $ws.Range("D22") '13.22'
Set-TextValue $ws.Range("E22") '  +0.63%  '

Set-TextValue $ws.Range("D23") '305.04'
Set-TextValue $ws.Range("E23") '  +1.98%  '

Set-TextValue $ws.Range("D24") '75.92'
Set-TextValue $ws.Range("E24") '  +0.79%  '

Set-TextValue $ws.Range("E25") '  -0.52%  '

Set-TextValue $ws.Range("D26") '28.81'
Set-TextValue $ws.Range("E26") '  +1.63%  '

Set-TextValue $ws.Range("D27") '4.49'
Set-TextValue $ws.Range("E27") '  +2.53%  '

Set-TextValue $ws.Range("B28") 'RenderToken'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D28") '7.73'
Set-TextValue $ws.Range("E28") '  +4.87%  '

Set-TextValue $ws.Range("B29") 'Filecoin'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D29") '7.94'
Set-TextValue $ws.Range("E29") '  -2.33%  '

Set-TextValue $ws.Range("D30") '0.176'
Set-TextValue $ws.Range("E30") '  +3.59%  '

Set-TextValue $ws.Range("E31") '  +5.72%  '

Set-TextValue $ws.Range("E32") '  -0.01%  '

Set-TextValue $ws.Range("D33") '11.57'
Set-TextValue $ws.Range("E33") '  +4.07%  '

Set-TextValue $ws.Range("D34") '40.31'
Set-TextValue $ws.Range("E34") '  +7.33%  '

Set-TextValue $ws.Range("D35") '0.0525'
Set-TextValue $ws.Range("E35") '  +7.49%  '

Set-TextValue $ws.Range("D36") '2.14'
Set-TextValue $ws.Range("E36") '  +0.03%  '

Set-TextValue $ws.Range("D37") '52.06'
Set-TextValue $ws.Range("E37") '  +0.28%  '

Set-TextValue $ws.Range("D38") '3.16'
Set-TextValue $ws.Range("E38") '  +0.59%  '

Set-TextValue $ws.Range("D39") '0.999'
Set-TextValue $ws.Range("E39") '  +0.03%  '

Set-TextValue $ws.Range("D40") '3.41'
Set-TextValue $ws.Range("E40") '  -3.32%  '

Set-TextValue $ws.Range("D41") '137.72'
Set-TextValue $ws.Range("E41") '  +3.22%  '

Set-TextValue $ws.Range("E42") '  +2.26%  '

Set-TextValue $ws.Range("E43") '  -0.10%  '

Set-TextValue $ws.Range("B44") 'Celestia'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D44") '17.28'
Set-TextValue $ws.Range("E44") '  -2.03%  '

Set-TextValue $ws.Range("B45") 'NEARProtocol'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D45") '4.00'
Set-TextValue $ws.Range("E45") '  +0.43%  '

Set-TextValue $ws.Range("D46") '0.282'
Set-TextValue $ws.Range("E46") '  -1.02%  '

Set-TextValue $ws.Range("B47") 'EnergySwap'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D47") '22.58'
Set-TextValue $ws.Range("E47") '  +1.31%  '

Set-TextValue $ws.Range("B48") 'WEMIXToken'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D48") '2.25'
Set-TextValue $ws.Range("E48") '  +7.50%  '

Set-TextValue $ws.Range("D49") '2.211.70'
Set-TextValue $ws.Range("E49") '  +2.82%  '

Set-TextValue $ws.Range("E50") '  -0.91%  '

Set-TextValue $ws.Range("E51") '  -11.46%  '
